# Update countries & provincias Spain
# Refresh COVID-19 stats for a set of countries, update the "last updated"
# timestamp, and re-sort the country table by total cases (column B) desc.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data range of countries (row 4 .. row 219), columns A (country) .. H (muertes)
$dataRange = $ws.Range("A4:H219")

# country -> updated stats: Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes
$updates = @{
    "Estados Unidos"       = @(2637180, 103,  1093527, 1415215, 0, 1,   128438)
    "India"                = @(549986,  789,  322161,  211321,  0, 17,  16504)
    "Iran"                 = @(225205,  2536, 186180,  28355,   0, 162, 10670)
    "Alemania"             = @(194900,  36,   178100,  7771,    0, 0,   9029)
    "Banglades"            = @(141801,  4014, 57780,   82238,   0, 45,  1783)
    "Catar"                = @(95106,   693,  80170,   14823,   0, 3,   113)
    "Indonesia"            = @(55092,   1082, 23800,   28487,   0, 51,  2805)
    "Oman"                 = @(39060,   910,  22422,   16469,   0, 6,   169)
    "Filipinas"            = @(36438,   983,  9956,    25227,   0, 11,  1255)
    "Polonia"              = @(34154,   247,  20897,   11813,   0, 6,   1444)
    "Suiza"                = @(31652,   35,   29100,   590,     0, 0,   1962)
    "Rumania"              = @(26582,   269,  18912,   6036,    0, 22,  1634)
    "Israel"               = @(23989,   234,  17114,   6556,    0, 1,   319)
    "Kazajistan"           = @(21327,   547,  13008,   8136,    0, 5,   183)
    "Austria"              = @(17723,   69,   16420,   600,     0, 1,   703)
    "Ghana"                = @(17351,   609,  12994,   4245,    0, 0,   112)
    "Nepal"                = @(13248,   476,  3134,    10085,   0, 1,   29)
    "Marruecos"            = @(12248,   196,  8790,    3234,    0, 3,   224)
    "Chequia"              = @(11604,   1,    7710,    3546,    0, 0,   348)
    "Malasia"              = @(8637,    3,    8334,    182,     0, 0,   121)
    "Finlandia"            = @(7209,    11,   6600,    281,     0, 0,   328)
    "Senegal"              = @(6698,    112,  4341,    2249,    0, 3,   108)
    "Bosnia y Herzegovina" = @(4325,    390,  2366,    1775,    0, 6,   184)
    "Albania"              = @(2466,    64,   1438,    970,     0, 3,   58)
    "Madagascar"           = @(2138,    60,   966,     1152,    0, 2,   20)
    "Hong Kong"            = @(1204,    4,    1105,    92,      0, 0,   7)
    "Uganda"               = @(870,     11,   808,     62,      0, 0,   0)
    "Malta"                = @(670,     0,    639,     22,      0, 0,   9)
    "Montenegro"           = @(498,     17,   315,     172,     0, 0,   11)
}

foreach ($country in $updates.Keys) {
    $vals = $updates[$country]
    $cell = $ws.Range("A4:A219").Find($country)
    $r = $cell.Row
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 6).Value = $vals[4]
    $ws.Cells.Item($r, 7).Value = $vals[5]
    $ws.Cells.Item($r, 8).Value = $vals[6]
}

# Re-sort the whole table by "Casos totales" (column B) descending, since the
# refreshed figures change the country ranking.
$dataRange.Sort($ws.Range("B4:B219"), 2)

# Update the "last updated" timestamp shown above the table.
$ws.Range("A1").Value = "Datos actualizados a 29 de Junio de 2020 a las 13:09"
